$wb = $excel.ActiveWorkbook

# --- Add the new "Czech" sheet as a copy of "Belgium" (same layout/styles), ---
# --- placed immediately after "Belgium" (i.e. as the last tab). ---
$belgium = $wb.Worksheets.Item("Belgium")
$belgium.Copy($null, $belgium)
$czech = $wb.Worksheets.Item($wb.Worksheets.Count)
$czech.Name = "Czech"

# Update the test-data reference number for Czech on the new sheet.
$czech.Range("B4").Value = "NGC-3477/T1731"

# Match the column widths recorded for the new sheet.
$czech.Columns.Item(2).ColumnWidth = 14.385416666666666
$czech.Columns.Item(3).ColumnWidth = 16.166666666666668
$czech.Columns.Item(4).ColumnWidth = 32.166666666666664

# Czech becomes the active sheet/tab, with D11 selected.
$czech.Range("D11").Select() | Out-Null
$czech.Activate() | Out-Null

# --- Belgium is no longer the selected tab; clear its old single-cell ---
# --- selection and select the full used range instead. ---
$belgium.Range("A1:D13").Select() | Out-Null

# Re-activate Czech so it ends up as the workbook's active tab/view.
$czech.Activate() | Out-Null
